$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I64").Value = 4827.3335
$ws.Range("H64").Value = 5296.4
$ws.Range("K64").Value = 4827.3335
$ws.Range("M64").Value = -4579.3335
$ws.Range("K67").Value = 4827.3335
$ws.Range("H67").Value = 5296.4
$ws.Range("I67").Value = 4827.3335
$ws.Range("M67").Value = -3969.3335
$ws.Range("H69").Value = 12290.714
$ws.Range("N70").Value = -10302
$ws.Range("J70").Value = 3254
$ws.Range("H70").Value = 2466.4814
$ws.Range("L70").Value = 9762
$ws.Range("H72").Value = 12290.714
$ws.Range("J73").Value = 3254
$ws.Range("N73").Value = -11634
$ws.Range("H73").Value = 2466.4814
$ws.Range("L73").Value = 9762
$ws.Range("I106").Value = 86195.336
$ws.Range("K106").Value = 86195.336
$ws.Range("M106").Value = -85564.336
$ws.Range("H106").Value = 65365.25
$ws.Range("I113").Value = 2500
$ws.Range("M113").Value = 754
$ws.Range("K113").Value = 2500
$ws.Range("H113").Value = 2500
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 100000
$ws.Range("H134").Value = 100000
$ws.Range("N134").Value = -110140
$ws.Range("M138").Value = 609.25
$ws.Range("K138").Value = 4530.75
$ws.Range("J138").Value = 3021.8
$ws.Range("I138").Value = 1510.25
$ws.Range("N138").Value = -19345.4
$ws.Range("H138").Value = 1954.8235
$ws.Range("L138").Value = 9065.400000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J56").Value = 0
$ws.Range("H56").Value = 1000000
$ws.Range("L56").Value = 0
$ws.Range("N104").Value = -24988
$ws.Range("H104").Value = 18000
$ws.Range("L104").Value = 18000
$ws.Range("J104").Value = 18000
$ws.Range("L105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("H105").Value = 0
$ws.Range("H107").Value = 70981.664
$ws.Range("N107").Value = -78661.664
$ws.Range("L107").Value = 70981.664
$ws.Range("J107").Value = 70981.664
$ws.Range("H109").Value = 69375.71000000001
$ws.Range("N109").Value = -72149.71000000001
$ws.Range("L109").Value = 69375.71000000001
$ws.Range("J109").Value = 69375.71000000001
$ws.Range("H111").Value = 69644
$ws.Range("N111").Value = -77824
$ws.Range("J111").Value = 69644
$ws.Range("L111").Value = 69644
$ws.Range("J114").Value = 45331.668
$ws.Range("H114").Value = 45331.668
$ws.Range("N114").Value = -54009.668
$ws.Range("L114").Value = 45331.668
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("H117").Value = 224373.33
$ws.Range("J117").Value = 224373.33
$ws.Range("N117").Value = -233551.33
$ws.Range("L117").Value = 224373.33
$ws.Range("N118").Value = -92487.836
$ws.Range("J118").Value = 89173.836
$ws.Range("H118").Value = 89173.836
$ws.Range("L118").Value = 89173.836
$ws.Range("H119").Value = 103212.125
$ws.Range("L119").Value = 103212.125
$ws.Range("J119").Value = 103212.125
$ws.Range("N119").Value = -112888.125
$ws.Range("L120").Value = 69000
$ws.Range("N120").Value = -78676
$ws.Range("H120").Value = 69000
$ws.Range("J120").Value = 69000
$ws.Range("J121").Value = 224379.17
$ws.Range("N121").Value = -227873.17
$ws.Range("L121").Value = 224379.17
$ws.Range("H121").Value = 224379.17
$ws.Range("N56").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("N115").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N31").Value = -2072.8
$ws.Range("H31").Value = 1596.0476
$ws.Range("L31").Value = 1482.8
$ws.Range("J31").Value = 1482.8
$ws.Range("H34").Value = 1596.0476
$ws.Range("N34").Value = -1886.8
$ws.Range("J34").Value = 1482.8
$ws.Range("L34").Value = 1482.8
$ws.Range("N75").Value = -35996
$ws.Range("L75").Value = 34000
$ws.Range("J75").Value = 34000
$ws.Range("H75").Value = 34000
$ws.Range("H78").Value = 34000
$ws.Range("N78").Value = -111984
$ws.Range("J78").Value = 34000
$ws.Range("L78").Value = 102000

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M2").Value = -15220
$ws.Range("J2").Value = 0
$ws.Range("H2").Value = 2555.5
$ws.Range("L2").Value = 0
$ws.Range("K2").Value = 15333
$ws.Range("I2").Value = 2555.5
$ws.Range("K12").Value = 45
$ws.Range("J12").Value = 548.9375
$ws.Range("L12").Value = 1646.8125
$ws.Range("M12").Value = 128
$ws.Range("N12").Value = -1992.8125
$ws.Range("H12").Value = 517.5294
$ws.Range("I12").Value = 15
$ws.Range("J132").Value = 3827.5557
$ws.Range("L132").Value = 34448.0013
$ws.Range("H132").Value = 2338.7827
$ws.Range("N132").Value = -39508.0013
$ws.Range("N2").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 1354.5555
$ws.Range("L80").Value = 2487.6667
$ws.Range("J80").Value = 2487.6667
$ws.Range("H80").Value = 1807.8
$ws.Range("M80").Value = -356.5554999999999
$ws.Range("N80").Value = -4483.6667
$ws.Range("K80").Value = 1354.5555
$ws.Range("H83").Value = 1807.8
$ws.Range("I83").Value = 1354.5555
$ws.Range("L83").Value = 12438.3335
$ws.Range("M83").Value = -1780.7775
$ws.Range("N83").Value = -22422.3335
$ws.Range("J83").Value = 2487.6667
$ws.Range("K83").Value = 6772.7775
$ws.Range("J102").Value = 11653.5
$ws.Range("M102").Value = -55773.74
$ws.Range("N102").Value = -14897.5
$ws.Range("H102").Value = 37989.938
$ws.Range("L102").Value = 11653.5
$ws.Range("I102").Value = 57395.74
$ws.Range("K102").Value = 57395.74
$ws.Range("J132").Value = 10999.25
$ws.Range("L132").Value = 32997.75
$ws.Range("H132").Value = 11249.375
$ws.Range("K132").Value = 34498.5
$ws.Range("M132").Value = -31968.5
$ws.Range("I132").Value = 11499.5
$ws.Range("N132").Value = -38057.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L6").Value = 42200
$ws.Range("N6").Value = -42424
$ws.Range("J6").Value = 42200
$ws.Range("H6").Value = 42200
$ws.Range("M40").Value = -5626.1113
$ws.Range("L40").Value = 0
$ws.Range("I40").Value = 5762.1113
$ws.Range("H40").Value = 5762.1113
$ws.Range("K40").Value = 5762.1113
$ws.Range("J40").Value = 0
$ws.Range("J61").Value = 494.5
$ws.Range("M61").Value = -25000368
$ws.Range("N61").Value = -898.5
$ws.Range("H61").Value = 22727836
$ws.Range("L61").Value = 494.5
$ws.Range("I61").Value = 25000570
$ws.Range("K61").Value = 25000570
$ws.Range("N68").Value = -8497
$ws.Range("M68").Value = -50
$ws.Range("L68").Value = 6999
$ws.Range("H68").Value = 4341.857
$ws.Range("K68").Value = 799
$ws.Range("J68").Value = 6999
$ws.Range("I68").Value = 799
$ws.Range("L71").Value = 34995
$ws.Range("J71").Value = 6999
$ws.Range("N71").Value = -42483
$ws.Range("I71").Value = 799
$ws.Range("M71").Value = -251
$ws.Range("K71").Value = 3995
$ws.Range("H71").Value = 4341.857
$ws.Range("N95").Value = -53581.082
$ws.Range("J95").Value = 48089.082
$ws.Range("H95").Value = 48089.082
$ws.Range("L95").Value = 48089.082
$ws.Range("M99").Value = -24004
$ws.Range("K99").Value = 26999
$ws.Range("I99").Value = 26999
$ws.Range("H99").Value = 38499.5
$ws.Range("I113").Value = 25000570
$ws.Range("M113").Value = -24998400
$ws.Range("L113").Value = 494.5
$ws.Range("K113").Value = 25000570
$ws.Range("N113").Value = -4834.5
$ws.Range("H113").Value = 22727836
$ws.Range("J113").Value = 494.5
$ws.Range("I136").Value = 2103.3635
$ws.Range("H136").Value = 38050.92
$ws.Range("M136").Value = -3760.0905
$ws.Range("K136").Value = 6310.0905
$ws.Range("N40").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I62").Value = 7412.75
$ws.Range("M62").Value = -6788.75
$ws.Range("K62").Value = 7412.75
$ws.Range("H62").Value = 9338.317999999999
$ws.Range("H65").Value = 9338.317999999999
$ws.Range("I65").Value = 7412.75
$ws.Range("M65").Value = -33943.75
$ws.Range("K65").Value = 37063.75
$ws.Range("H107").Value = 748.5
$ws.Range("I107").Value = 0
$ws.Range("N107").Value = -6085.5
$ws.Range("L107").Value = 2245.5
$ws.Range("K107").Value = 0
$ws.Range("J107").Value = 748.5
$ws.Range("M107").ClearContents()
